$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: "Class 4" summary cell becomes "Total SDS / SHOP" ---
$ws.Range("G9").Value = "Total SDS / SHOP"

# --- Row 11: Remark header (text unchanged, index re-sequenced internally) ---
$ws.Range("A11").Value = "Remark : Classification of Hazardous Substances (ประเภทของวัตถุอันตราย)"

# --- Prime formatting for the cells that are new / previously unformatted ---
# Column A style (from existing "Class 1" cell) -> A16 (was blank) and new rows A17:A22
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A17:A22").PasteSpecial(-4122) | Out-Null
# Column C style -> new rows C17:C22
$ws.Range("C12").Copy() | Out-Null
$ws.Range("C17:C22").PasteSpecial(-4122) | Out-Null
# Column G style -> new rows G17:G22
$ws.Range("G12").Copy() | Out-Null
$ws.Range("G17:G22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 12 ---
$ws.Range("A12").Value = "Class 1"
$ws.Range("C12").Value = "วัตถุระเบิด"
$ws.Range("G12").Value = "Class 6.1A สารติดไฟได้ ที่มีคุณสมบัติเป็นพิษ"

# --- Row 13 ---
$ws.Range("A13").Value = "Class 2A"
$ws.Range("C13").Value = "ก๊าซอัด ก๊าซเหลว ก๊าซละลายได้ภายใต้ความดัน"
$ws.Range("G13").Value = "Class 6.1B สารไม่ติดไฟ ที่มีคุณสมบัติเป็นพิษ"

# --- Row 14 ---
$ws.Range("A14").Value = "Class 2B"
$ws.Range("C14").Value = "ก๊าซภายใต้ความดันในกระป๋องเสปร์ย"
$ws.Range("G14").Value = "Class 6.2 สารติดเชื้อ"

# --- Row 15 ---
$ws.Range("A15").Value = "Class 3A"
$ws.Range("C15").Value = "ของเหลวไวไฟ จุดวาบไฟไม่เกิน 60 ºC"
$ws.Range("G15").Value = "Class 7 สารกัมมันตรังสี"

# --- Row 16 ---
$ws.Range("A16").Value = "Class 3B"
$ws.Range("C16").Value = "ของเหลวไวไฟ จุดวาบไฟมากกว่า 60 ºC-93 ºC คุณสมบัติเข้ากับน้ำไม่ได้"
$ws.Range("G16").Value = "Class 8A สารติดไฟ ที่มีคุณสมบัติกัดกร่อน"

# --- Row 17 (new) ---
$ws.Range("A17").Value = "Class 4.1A"
$ws.Range("C17").Value = "ของแข็งไวไฟ ที่มีคุณสมบัติระเบิด"
$ws.Range("G17").Value = "Class 8B สารไม่ติดไฟ ที่มีคุณสมบัติกัดกร่อน"

# --- Row 18 (new) ---
$ws.Range("A18").Value = "Class 4.1B"
$ws.Range("C18").Value = "ของแข็งไวไฟ ที่ไม่มีคุณสมบัติระเบิด"
$ws.Range("G18").Value = "Class 9 วัตถุอันตรายประเภทอื่นๆ"

# --- Row 19 (new) ---
$ws.Range("A19").Value = "Class 4.2"
$ws.Range("C19").Value = "สารที่มีความเสี่ยงต่อการลุกไหม้ได้เอง"
$ws.Range("G19").Value = "Class 10 ของเหลวติดไฟ"

# --- Row 20 (new) ---
$ws.Range("A20").Value = "Class 4.3"
$ws.Range("C20").Value = "สารให้ก๊าซไวไฟ เมื่อสัมผัสกับน้ำ"
$ws.Range("G20").Value = "Class 11 ของแข็งติดไฟได้"

# --- Row 21 (new) ---
$ws.Range("A21").Value = "Class 5.1"
$ws.Range("C21").Value = "สารออกซิไดซ์"
$ws.Range("G21").Value = "Class 12 ของเหลวไม่ติดไฟ"

# --- Row 22 (new) ---
$ws.Range("A22").Value = "Class 5.2"
$ws.Range("C22").Value = "สารเปอร์ออกซิไดซ์"
$ws.Range("G22").Value = "Class 13 ของแข็งไม่ติดไฟ"

# --- Row heights for row 16 (now populated) and the newly added rows 17:22 ---
$ws.Range("16:22").RowHeight = 23.25

# --- Match the active selection left by the author ---
$ws.Range("E12").Select() | Out-Null
